# IndexTOC.docx edit: expand several bullet labels in the "Models" /
# "Super / Sub Context Alignment" section with attribute/object wording,
# and append a sentence to the long "Context: State factor ..." paragraph.
#
# Find.Execute params (by position):
#   1 FindText, 2 MatchCase, 3 MatchWholeWord, 4 MatchWildcards,
#   5 MatchSoundsLike, 6 MatchAllWordForms, 7 Forward, 8 Wrap,
#   9 Format, 10 ReplaceWith, 11 Replace (wdReplaceAll = 2)

$d = $word.ActiveDocument

$d.Content.Find.Execute("Reified Statements:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Reified Statements (attributes), SPO Factors x Kinds (objects):", 2)

$d.Content.Find.Execute("Graph Statements (materialize / align):", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Graph Statements (materialize / align) URNFactors attributes. Objects: CSPO Factor:", 2)

$d.Content.Find.Execute("SubjectKind:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "SubjectKind attributes. Objects: aggregated Subject Factors:", 2)

$d.Content.Find.Execute("PredicateKind:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "PredicateKind attributes. Objects: aggregated Predicate Factors.", 2)

$d.Content.Find.Execute("ObjectKind:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ObjectKind attributes. Objects: aggregated Object Factors.", 2)

$d.Content.Find.Execute("Order / Flow Mappings:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Order / Flow Mappings attributes. Objects: aggregated State Factors.", 2)

$d.Content.Find.Execute("Joins matching Kinds factors (flows).", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Joins matching Kinds factors (flows). Provenance state flows.", 2)
